$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PercentText($addr, $val, $fmtSrc) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $ws.Range($fmtSrc).Copy()
    $c.PasteSpecial(-4122)
}

$ws.Range('E2').Value = '2026-02-10 06:48:18'
$ws.Range('E3').Value = '2026-02-10 06:48:20'
$ws.Range('G3').Value = '185 cm'
$ws.Range('I3').Value = '8.2 mm'
$ws.Range('E4').Value = '2026-02-10 06:48:22'
$ws.Range('E5').Value = '2026-02-10 06:48:25'
$ws.Range('G5').Value = '132 cm'
$ws.Range('I5').Value = '13.2 mm'
$ws.Range('E6').Value = '2026-02-10 06:48:27'
$ws.Range('N6').Value = '6.3 °C 6:29 TU'
$ws.Range('E7').Value = '2026-02-10 06:48:29'
$ws.Range('E8').Value = '2026-02-10 06:48:32'
Set-PercentText 'H8' '94%' 'F8'
$ws.Range('N8').Value = '8.2 °C 6:00 TU'
$ws.Range('O8').Value = '8.5 °C'
$ws.Range('E9').Value = '2026-02-10 06:48:34'
$ws.Range('N9').Value = '4.9 °C 6:14 TU'
$ws.Range('O9').Value = '6.4 °C'
$ws.Range('E10').Value = '2026-02-10 06:48:36'
$ws.Range('N10').Value = '4.5 °C 6:01 TU'
$ws.Range('O10').Value = '6.9 °C'
$ws.Range('E11').Value = '2026-02-10 06:48:39'
$ws.Range('E12').Value = '2026-02-10 06:48:41'
$ws.Range('O12').Value = '6.7 °C'
$ws.Range('E13').Value = '2026-02-10 06:48:43'
$ws.Range('I13').Value = '2.3 mm'
$ws.Range('J13').Value = '1008.3 hPa'
$ws.Range('N13').Value = '2.4 °C 6:02 TU'
$ws.Range('E14').Value = '2026-02-10 06:48:46'
$ws.Range('E15').Value = '2026-02-10 06:48:48'
$ws.Range('O15').Value = '6.2 °C'
$ws.Range('E16').Value = '2026-02-10 06:48:50'
$ws.Range('G16').Value = '79 cm'
Set-PercentText 'H16' '89%' 'F16'
$ws.Range('I16').Value = '12.8 mm'
$ws.Range('E17').Value = '2026-02-10 06:48:52'
Set-PercentText 'H17' '90%' 'F17'
$ws.Range('O17').Value = '2.9 °C'
$ws.Range('E18').Value = '2026-02-10 06:48:55'
$ws.Range('J18').Value = '1005.2 hPa'
$ws.Range('N18').Value = '4.5 °C 6:25 TU'
$ws.Range('O18').Value = '7.1 °C'
$ws.Range('E19').Value = '2026-02-10 06:48:57'
$ws.Range('N19').Value = '3.2 °C 6:29 TU'
$ws.Range('E20').Value = '2026-02-10 06:48:59'
$ws.Range('I20').Value = '2.2 mm'
$ws.Range('M20').Value = '-0.4 °C 6:19 TU'
$ws.Range('E21').Value = '2026-02-10 06:49:01'
$ws.Range('I21').Value = '3.4 mm'
$ws.Range('J21').Value = '1007.6 hPa'
$ws.Range('E22').Value = '2026-02-10 06:49:04'
$ws.Range('M22').Value = '-1.3 °C 6:18 TU'
$ws.Range('E23').Value = '2026-02-10 06:49:06'
$ws.Range('I23').Value = '11.2 mm'
$ws.Range('E24').Value = '2026-02-10 06:49:08'
$ws.Range('I24').Value = '1.6 mm'
$ws.Range('N24').Value = '8.0 °C 6:29 TU'
$ws.Range('E25').Value = '2026-02-10 06:49:11'
$ws.Range('I25').Value = '7.1 mm'
$ws.Range('O25').Value = '-0.6 °C'
$ws.Range('E26').Value = '2026-02-10 06:49:13'
Set-PercentText 'H26' '86%' 'F26'
$ws.Range('J26').Value = '1005.0 hPa'
$ws.Range('M26').Value = '5.9 °C 6:29 TU'
$ws.Range('O26').Value = '3.2 °C'
$ws.Range('E27').Value = '2026-02-10 06:49:16'
$ws.Range('G27').Value = '173 cm'
$ws.Range('I27').Value = '2.0 mm'
$ws.Range('E28').Value = '2026-02-10 06:49:18'
$ws.Range('N28').Value = '3.3 °C 6:00 TU'
$ws.Range('O28').Value = '4.9 °C'
$ws.Range('E29').Value = '2026-02-10 06:49:20'
$ws.Range('E30').Value = '2026-02-10 06:49:23'
$ws.Range('O30').Value = '7.3 °C'
$ws.Range('E31').Value = '2026-02-10 06:49:25'
$ws.Range('E32').Value = '2026-02-10 06:49:27'
$ws.Range('E33').Value = '2026-02-10 06:49:30'
$ws.Range('I33').Value = '5.3 mm'
$ws.Range('J33').Value = '1007.9 hPa'
$ws.Range('E34').Value = '2026-02-10 06:49:32'
Set-PercentText 'H34' '81%' 'F34'
$ws.Range('I34').Value = '2.8 mm'
$ws.Range('O34').Value = '2.3 °C'
$ws.Range('E35').Value = '2026-02-10 06:49:35'
$ws.Range('M35').Value = '11.0 °C 6:27 TU'
$ws.Range('E36').Value = '2026-02-10 06:49:37'
$ws.Range('N36').Value = '6.2 °C 6:14 TU'
$ws.Range('O36').Value = '8.6 °C'
$ws.Range('E37').Value = '2026-02-10 06:49:40'
$ws.Range('O37').Value = '3.7 °C'
$ws.Range('E38').Value = '2026-02-10 06:49:42'
$ws.Range('E39').Value = '2026-02-10 06:49:44'
$ws.Range('I39').Value = '3.0 mm'
$ws.Range('E40').Value = '2026-02-10 06:49:47'
$ws.Range('I40').Value = '4.0 mm'
$ws.Range('J40').Value = '1008.4 hPa'
$ws.Range('O40').Value = '4.6 °C'
$ws.Range('E41').Value = '2026-02-10 06:49:49'
$ws.Range('J41').Value = '1005.1 hPa'
$ws.Range('O41').Value = '9.9 °C'
$ws.Range('E42').Value = '2026-02-10 06:49:52'
$ws.Range('N42').Value = '6.5 °C 6:08 TU'
$ws.Range('O42').Value = '7.9 °C'
$ws.Range('E43').Value = '2026-02-10 06:49:54'
$ws.Range('N43').Value = '5.4 °C 6:29 TU'
$ws.Range('O43').Value = '6.0 °C'
$ws.Range('E44').Value = '2026-02-10 06:49:56'
$ws.Range('I44').Value = '7.6 mm'
$ws.Range('E45').Value = '2026-02-10 06:49:59'
$ws.Range('I45').Value = '17.2 mm'
$ws.Range('M45').Value = '4.4 °C 6:03 TU'
$ws.Range('O45').Value = '3.3 °C'
$ws.Range('E46').Value = '2026-02-10 06:50:01'
Set-PercentText 'H29' '93%' 'F29'
$ws.Range('I29').Value = '0.6 mm'
$ws.Range('K29').Value = '0.0 MJ/m2'
$ws.Range('L29').Value = '16.2 km/h - 211º 0:18 TU'
$ws.Range('M29').Value = '10.6 °C 2:34 TU'
$ws.Range('N29').Value = '5.5 °C 6:16 TU'
$ws.Range('O29').Value = '8.5 °C'

$excel.CutCopyMode = $false
